# Build site at 2023-04-12 14:53:07 UTC
# LOM3258.xlsx update:
#  - Objectives (PT) text filled in (row 10)
#  - A new row is inserted (blank label) holding the professor's name under
#    "Docentes responsaveis:" (row 12), which previously was mis-placed two
#    rows further down under "Metodo:"
#  - "Programa resumido:" / "Programa:" (PT short-syllabus / syllabus) text
#    filled in
#  - A new "Bibliografia:" row with real bibliography text is appended
#  - Column A's width definition is narrowed back to just column A

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fill in the (previously empty) Portuguese "Objetivos:" text - row 10
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = 'Proporcionar ao aluno ingressante de Engenharia Física os conhecimentos práticos de eletrônica e computação física com microcontrolador Arduino visando sua aplicação em projetos científicos e tecnológicos.'
$ws.Range("C10").Value = 'Proporcionar ao aluno ingressante de Engenharia Física os conhecimentos práticos de eletrônica e computação física com microcontrolador Arduino visando sua aplicação em projetos científicos e tecnológicos.'

# ---------------------------------------------------------------------------
# 2) Insert a new row at 13 - this shifts the old rows 13..21 down to 14..22,
#    fixing the label (column A) / value (columns B & C) misalignment that
#    existed from row 13 onward.
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Insert()

# The freshly inserted row 13 only has an empty, styled A13 placeholder cell.
# Give B13/C13 the same look (style) as the row right below (old row 13,
# "Programa resumido" / "Semestral", now shifted to row 14), then clear A13
# and set the real value: the professor's name belongs here, right under
# "Docentes responsaveis:" (row 12).
$ws.Range("B14").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$ws.Range("A13").Clear()
$ws.Range("B13").Value = '519033 - Carlos Yujiro Shigue'
$ws.Range("C13").Value = '519033 - Carlos Yujiro Shigue'

# ---------------------------------------------------------------------------
# 3) Fill in the Portuguese short syllabus text - row 14 ("Programa resumido:")
# ---------------------------------------------------------------------------
$ws.Range("B14").Value = 'Introdução ao Arduino. Conceitos de eletrônica analógica e digital. Montagem de circuitos eletrônicos básicos. Programação e controle de circuitos eletrônicos em linguagem C. Aplicação e desenvolvimento de projetos baseados em Arduino.'
$ws.Range("C14").Value = 'Introdução ao Arduino. Conceitos de eletrônica analógica e digital. Montagem de circuitos eletrônicos básicos. Programação e controle de circuitos eletrônicos em linguagem C. Aplicação e desenvolvimento de projetos baseados em Arduino.'

# ---------------------------------------------------------------------------
# 4) Fill in the Portuguese full syllabus text - row 16 ("Programa:")
# ---------------------------------------------------------------------------
$ws.Range("B16").Value = 'Introdução ao microcontrolador Arduino: histórico, tipos e recursos. Oficina prática: instalação e configuração do IDE Arduino.Conceitos básicos de eletrônica: funcionamento da protoboard, componentes e instrumentos eletrônicos, medições com multímetro e osciloscópio. Grandezas elétricas: resistência, tensão e corrente. Oficina: montagem de circuitos eletrônicos.Introdução à linguagem de programação Wiring baseada em C/C++. Tipos de dados, sintaxe básica, controle de fluxo, funções da biblioteca padrão. Principais bibliotecasEntradas e saídas do Arduino. Sinais analógicos e digitais.Controle de dispositivos utilizando PWM.Eletrônica analógica. Conversores analógico-digitais do Arduino. Oficina: leitura de dados de sensores. Comunicação serial/USB com o PC. Utilização do Monitor Serial da IDE.Controle de motor cc e servomotor com PWM. Controle de potência com relé e SSR.Tópicos avançados: comunicação Ethernet com Arduino. Comunicação sem fio via Bluetooth.Armazenamento de dados utilizando a EEPROM do ATMega328 e cartão de memória SD.Desenvolvimento de software de qualidade.Desenvolvimento de projetos utilizando microcontrolador Arduino.'
$ws.Range("C16").Value = 'Introdução ao microcontrolador Arduino: histórico, tipos e recursos. Oficina prática: instalação e configuração do IDE Arduino.Conceitos básicos de eletrônica: funcionamento da protoboard, componentes e instrumentos eletrônicos, medições com multímetro e osciloscópio. Grandezas elétricas: resistência, tensão e corrente. Oficina: montagem de circuitos eletrônicos.Introdução à linguagem de programação Wiring baseada em C/C++. Tipos de dados, sintaxe básica, controle de fluxo, funções da biblioteca padrão. Principais bibliotecasEntradas e saídas do Arduino. Sinais analógicos e digitais.Controle de dispositivos utilizando PWM.Eletrônica analógica. Conversores analógico-digitais do Arduino. Oficina: leitura de dados de sensores. Comunicação serial/USB com o PC. Utilização do Monitor Serial da IDE.Controle de motor cc e servomotor com PWM. Controle de potência com relé e SSR.Tópicos avançados: comunicação Ethernet com Arduino. Comunicação sem fio via Bluetooth.Armazenamento de dados utilizando a EEPROM do ATMega328 e cartão de memória SD.Desenvolvimento de software de qualidade.Desenvolvimento de projetos utilizando microcontrolador Arduino.'

# ---------------------------------------------------------------------------
# 5) Rows 19/20/21 ("Metodo:" / "Criterio:" / "Norma de recuperacao:") now
#    line up correctly with the post-shift label thanks to the insert above,
#    so no value changes are required there.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 6) Row 22 ("Bibliografia:") - replace the mis-placed recovery-policy text
#    (which the shift already moved up to row 21, where it belongs) with the
#    real bibliography text.
# ---------------------------------------------------------------------------
$ws.Range("B22").Value = "BANZI, M. Primeiros passos com o Arduino, São Paulo: O´Reilly Novatec, 2010.`nMcROBERTS, M. Arduino Básico, São Paulo: Novatec, 2011.`nMONK, S. Programação com Arduino, Porto Alegre: Bookman Editora, 2013.`nMONK, S. Programação com Arduino II, Porto Alegre: Bookman Editora, 2015.`nBLUM, J. Exploring Arduino, New York: John Wiley, 2013."
$ws.Range("C22").Value = "BANZI, M. Primeiros passos com o Arduino, São Paulo: O´Reilly Novatec, 2010.`nMcROBERTS, M. Arduino Básico, São Paulo: Novatec, 2011.`nMONK, S. Programação com Arduino, Porto Alegre: Bookman Editora, 2013.`nMONK, S. Programação com Arduino II, Porto Alegre: Bookman Editora, 2015.`nBLUM, J. Exploring Arduino, New York: John Wiley, 2013."

# ---------------------------------------------------------------------------
# 7) Column A's width definition previously (incorrectly) spanned columns
#    A:B; narrow it back down to just column A (column B keeps its own,
#    wider, definition right below).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 29.877604166666668
